# Generate Report for Handback
# Refresh the handback timestamps recorded on the Overview, zh-cn and de-de
# sheets to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the zh-cn handback row
$wsOverview.Range("G2").Value = "2016-08-18 08:50:35"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first (7c794a4b...) row
$wsZhCn.Range("H2").Value = "2016-08-18 08:50:30"
$wsZhCn.Range("K2").Value = "2016-08-18 08:50:46"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first (7c794a4b...) row
$wsDeDe.Range("H2").Value = "2016-08-18 08:50:35"
$wsDeDe.Range("K2").Value = "2016-08-18 08:50:55"
